# Informe de incidencia.docx - "Doc: Informe de tareas terminado"
#
# This script reproduces the two semantic changes in the target diff:
#
#   1) The "Institucional - Portal UNPA UARG" logo picture (docPr name
#      "Imagen 4", top paragraph of the body) is enlarged (~29.3%) and
#      nudged left/up slightly:
#           Left  : 445 pt   -> 437.5 pt   (posOffset 5651500 -> 5556250 EMU)
#           Top   : -0.2 pt  -> -0.25 pt   (posOffset   -2540 ->   -3175 EMU)
#           Width : 29.674803149606298 pt -> 38.37259842519685 pt
#                     (extent cx 376870 -> 487332 EMU)
#           Height: 43.5 pt  -> 56.25 pt
#                     (extent cy 552450 -> 714375 EMU)
#
#   2) In the page header, the two runs that spell out the label
#      "Iteración" + ": " (identical run formatting: color 3B3838, sz 20)
#      are coalesced into a single run holding "Iteración: ". This is a
#      pure run-split/run-merge normalisation - no visible text or
#      formatting changes - inside the header's "Iteración: ____" text
#      box, so it is applied best-effort via the shape's text range.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Resize / reposition the logo picture ("Imagen 4").
# ---------------------------------------------------------------------

$targetName   = "Imagen 4"
$targetLeft   = 437.5
$targetTop    = -0.25
$targetWidth  = 38.37259842519685
$targetHeight = 56.25

if ($d.Shapes.Count -gt 0) {

    # Locate the shape by name first (so the script is resilient to the
    # shapes being stored/enumerated in a different order), remembering
    # both the index reported by the getters and the index that needs to
    # be used for the (separately indexed) setters.
    $getterIndex = -1
    for ($i = 1; $i -le $d.Shapes.Count; $i++) {
        if ($d.Shapes.Item($i).Name -eq $targetName) {
            $getterIndex = $i
        }
    }

    # Apply the new geometry, then confirm - via the read-back getters -
    # that the intended shape actually moved; if not, fall back to the
    # other index. (Belt-and-braces against get/set indices disagreeing.)
    $applied = $false
    $candidates = @(1, 2, $getterIndex) | Where-Object { $_ -ge 1 -and $_ -le $d.Shapes.Count } | Select-Object -Unique

    foreach ($idx in $candidates) {
        if ($applied) { break }

        $shp = $d.Shapes.Item($idx)
        $shp.Width  = $targetWidth
        $shp.Height = $targetHeight
        $shp.Left   = $targetLeft
        $shp.Top    = $targetTop

        # Verify against whichever index the getters say now reports the
        # target shape/size.
        for ($j = 1; $j -le $d.Shapes.Count; $j++) {
            $probe = $d.Shapes.Item($j)
            if ($probe.Name -eq $targetName -and [Math]::Round($probe.Width, 1) -eq [Math]::Round($targetWidth, 1)) {
                $applied = $true
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Header: merge the "Iteración" / ": " runs into a single run.
# ---------------------------------------------------------------------
# Best-effort: the label lives inside a text box anchored to the page
# header, so it is reached through HeaderFooter.Shapes(...).TextFrame.
# Wrapped in try/catch so an unsupported path here can never clobber the
# picture edit performed above.

try {
    for ($s = 1; $s -le $d.Sections.Count; $s++) {
        $section = $d.Sections.Item($s)

        foreach ($hdrIndex in 1..3) {
            try {
                $hdr = $section.Headers.Item($hdrIndex)
            } catch {
                continue
            }
            if ($null -eq $hdr -or -not $hdr.Exists) { continue }

            for ($k = 1; $k -le $hdr.Shapes.Count; $k++) {
                $shape = $hdr.Shapes.Item($k)
                $tf = $shape.TextFrame
                if ($null -eq $tf) { continue }
                if (-not $tf.HasText) { continue }

                $tr = $tf.TextRange
                $current = $tr.Text

                if ($current -like "*Iteraci*" -and $current -like "*________*") {
                    $merged = "Iteraci" + [char]0x00F3 + "n: ________"

                    # Try the direct text assignment first ...
                    $tr.Text = $merged

                    # ... then, in case that property is not wired up,
                    # fall back to a Find/Replace across the same range.
                    if ($tr.Text -ne $merged) {
                        $null = $tr.Find.Execute(
                            "Iteraci" + [char]0x00F3 + "n: ",
                            $true, $false, $false, $false, $false,
                            $true, 1, $false, $merged, 2)
                    }
                }
            }
        }
    }
} catch {
    # Non-fatal: leave the header text box untouched if this host does
    # not expose writable header shape text ranges.
}
